$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking strings (e.g. "556.42", "65.087.43")
# that must stay as literal text, matching the source inlineStr cells.
# Force text formatting, assign, then restore the default "Normal" style
# so the cell keeps no explicit style (s attribute), exactly like the original.

$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '65.087.43'
$cell.Style = 'Normal'
$ws.Range('E2').Value = '  +0.75%  '
$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '3.374.09'
$cell.Style = 'Normal'
$ws.Range('E3').Value = '  +0.34%  '
$cell = $ws.Range('D4')
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$cell.Style = 'Normal'
$ws.Range('E4').Value = '  -0.03%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '556.42'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  -0.20%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '174.87'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  -0.62%  '
$ws.Range('E7').Value = '  +1.80%  '
$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '3.363.10'
$cell.Style = 'Normal'
$ws.Range('E8').Value = '  +0.22%  '
$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$cell.Style = 'Normal'
$ws.Range('E9').Value = '  +0.06%  '
$ws.Range('E10').Value = '  +5.91%  '
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '0.637'
$cell.Style = 'Normal'
$ws.Range('E11').Value = '  +1.04%  '
$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '53.71'
$cell.Style = 'Normal'
$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '0.0000279'
$cell.Style = 'Normal'
$ws.Range('E13').Value = '  +2.08%  '
$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '9.21'
$cell.Style = 'Normal'
$ws.Range('E14').Value = '  +1.44%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '3.907.19'
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  +0.41%  '
$ws.Range('E16').Value = '  -0.82%  '
$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '3.372.28'
$cell.Style = 'Normal'
$ws.Range('E17').Value = '  +0.43%  '
$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '0.118'
$cell.Style = 'Normal'
$ws.Range('E18').Value = '  -0.59%  '
$ws.Range('E19').Value = '  +0.21%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '64.837.37'
$cell.Style = 'Normal'
$ws.Range('E20').Value = '  +0.64%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  +1.40%  '
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '454.67'
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  -0.78%  '
$ws.Range('E23').Value = '  +1.51%  '
$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '4.08'
$cell.Style = 'Normal'
$ws.Range('E24').Value = '  -0.49%  '
$ws.Range('B25').Value = 'InternetComputer(DFINITY)'
$ws.Range('C25').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '14.08'
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  +5.35%  '
$ws.Range('B26').Value = 'Litecoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '87.53'
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  +2.43%  '
$ws.Range('E27').Value = '  +0.85%  '
$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '10.71'
$cell.Style = 'Normal'
$ws.Range('E28').Value = '  -0.91%  '
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '8.72'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  -0.96%  '
$ws.Range('E30').Value = '  +4.13%  '
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '6.54'
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  -0.98%  '
$ws.Range('B32').Value = 'OKB'
$ws.Range('C32').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '63.08'
$cell.Style = 'Normal'
$ws.Range('E32').Value = '  +7.58%  '
$ws.Range('B33').Value = 'Cosmos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '11.46'
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  -0.53%  '
$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '577.63'
$cell.Style = 'Normal'
$ws.Range('E34').Value = '  -0.66%  '
$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '0.108'
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  -0.40%  '
$ws.Range('E36').Value = '  +0.09%  '
$ws.Range('E37').Value = '  +4.27%  '
$ws.Range('E38').Value = '  +0.45%  '
$ws.Range('E39').Value = '  -0.24%  '
$ws.Range('E40').Value = '  +0.83%  '
$ws.Range('E41').Value = '  -1.93%  '
$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '3.092.00'
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  -0.34%  '
$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '0.0417'
$cell.Style = 'Normal'
$ws.Range('E43').Value = '  +1.51%  '
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '2.77'
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  -1.30%  '
$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '2.46'
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  -0.75%  '
$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '0.134'
$cell.Style = 'Normal'
$ws.Range('E46').Value = '  +2.23%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '3.18'
$cell.Style = 'Normal'
$ws.Range('E47').Value = '  -1.40%  '
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '142.40'
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  +5.52%  '
$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '0.998'
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  +0.03%  '
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '2.54'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  -1.76%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '8.28'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  -0.20%  '
